# Update the "取得日時" (acquired timestamp) column on the "ランサーズ" sheet.
# All data rows (2-10) previously stamped at 2025-09-27 12:32:16 are
# re-stamped to reflect the latest append run at 2025-09-27 12:41:01.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-27 12:41:01"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 10
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq "2025-09-27 12:32:16") {
        $cell.Value = $newTimestamp
    }
}

$wb.Save()
